$wb = $excel.ActiveWorkbook
$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-03-01 09:18:49"

$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-03-01 09:18:59"
